$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This workbook holds NATMI ligand-receptor pair output for Rspo1 -> Lgr5.
# The commit "update scripts wuth new tpm" re-ran the analysis with updated
# TPM values. The sending cluster changes from "Resolving-Mac" to
# "Inflammatory-Mac", the set of target clusters is re-ordered (row 4 becomes
# MuSCs, row 5 becomes Resolving-Mac, which is now also a valid target
# cluster), and all of the derived expression / specificity statistics are
# recalculated accordingly.
# ---------------------------------------------------------------------------

# Row 2 (target cluster: ECs)
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Rspo1"
$ws.Range("C2").Value = "Lgr5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7878926666666667
$ws.Range("H2").Value = 2.363678
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.4220680000000001
$ws.Range("N2").Value = 1.266204
$ws.Range("O2").Value = 0.06297272914451456
$ws.Range("P2").Value = 0.06297272914451456
$ws.Range("Q2").Value = 0.3325442820346667
$ws.Range("R2").Value = 2.992898538312001
$ws.Range("S2").Value = 0.06297272914451456
$ws.Range("T2").Value = 0.06297272914451456

# Row 3 (target cluster: FAPs)
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Rspo1"
$ws.Range("C3").Value = "Lgr5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7878926666666667
$ws.Range("H3").Value = 2.363678
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.445601666666667
$ws.Range("N3").Value = 7.336805
$ws.Range("O3").Value = 0.3648848321843242
$ws.Range("P3").Value = 0.3648848321843242
$ws.Range("Q3").Value = 1.926871618754445
$ws.Range("R3").Value = 17.34184456879
$ws.Range("S3").Value = 0.3648848321843242
$ws.Range("T3").Value = 0.3648848321843242

# Row 4 (target cluster: MuSCs)
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Rspo1"
$ws.Range("C4").Value = "Lgr5"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7878926666666667
$ws.Range("H4").Value = 2.363678
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.729698666666666
$ws.Range("N4").Value = 11.189096
$ws.Range("O4").Value = 0.5564726629989883
$ws.Range("P4").Value = 0.5564726629989883
$ws.Range("Q4").Value = 2.938602228343111
$ws.Range("R4").Value = 26.447420055088
$ws.Range("S4").Value = 0.5564726629989883
$ws.Range("T4").Value = 0.5564726629989883

# Row 5 (target cluster: Resolving-Mac)
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Rspo1"
$ws.Range("C5").Value = "Lgr5"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7878926666666667
$ws.Range("H5").Value = 2.363678
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.105025
$ws.Range("N5").Value = 0.315075
$ws.Range("O5").Value = 0.01566977567217282
$ws.Range("P5").Value = 0.01566977567217282
$ws.Range("Q5").Value = 0.08274842731666666
$ws.Range("R5").Value = 0.7447358458500001
$ws.Range("S5").Value = 0.01566977567217282
$ws.Range("T5").Value = 0.01566977567217282
